$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "Bar"
$ws.Range("E2").Value = "Bar"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "Bar"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = "-"
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "Bar"
$ws.Range("G3").Value = "Bar"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = "Bar"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "-"
$ws.Range("L3").Value = "-"
$ws.Range("M3").Value = "Bar"
$ws.Range("N3").Value = "Bar"
$ws.Range("O3").Value = "-"
$ws.Range("B4").Value = "Bar"
$ws.Range("C4").Value = "Bar"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "Bar"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "-"
$ws.Range("H4").Value = "-"
$ws.Range("I4").Value = "-"
$ws.Range("J4").Value = "-"
$ws.Range("K4").Value = "-"
$ws.Range("L4").Value = "-"
$ws.Range("M4").Value = "-"
$ws.Range("N4").Value = "-"
$ws.Range("O4").Value = "-"
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "-"
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "-"
$ws.Range("H5").Value = "-"
$ws.Range("I5").Value = "-"
$ws.Range("J5").Value = "Bar"
$ws.Range("K5").Value = "-"
$ws.Range("M5").Value = "-"
$ws.Range("O5").Value = "-"
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "-"
$ws.Range("G6").Value = "Bar"
$ws.Range("H6").Value = "-"
$ws.Range("J6").Value = "-"
$ws.Range("L6").Value = "-"
$ws.Range("M6").Value = "-"
$ws.Range("N6").Value = "-"
$ws.Range("O6").Value = "Bar"
$ws.Range("B7").Value = "-"
$ws.Range("D7").Value = "Server"
$ws.Range("F7").Value = "Server"
$ws.Range("G7").Value = "Server"
$ws.Range("H7").Value = "-"
$ws.Range("I7").Value = "Expo"
$ws.Range("J7").Value = "-"
$ws.Range("K7").Value = "Bar"
$ws.Range("L7").Value = "Server"
$ws.Range("O7").Value = "Bar"
$ws.Range("B8").Value = "Server"
$ws.Range("C8").Value = "Server"
$ws.Range("G8").Value = "Expo"
$ws.Range("H8").Value = "-"
$ws.Range("I8").Value = "Server"
$ws.Range("J8").Value = "-"
$ws.Range("K8").Value = "-"
$ws.Range("L8").Value = "-"
$ws.Range("B9").Value = "-"
$ws.Range("C9").Value = "-"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "Server"
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "-"
$ws.Range("H9").Value = "-"
$ws.Range("I9").Value = "-"
$ws.Range("J9").Value = "-"
$ws.Range("K9").Value = "-"
$ws.Range("L9").Value = "-"
$ws.Range("M9").Value = "-"
$ws.Range("N9").Value = "Server"
$ws.Range("O9").Value = "Server"
$ws.Range("B10").Value = "Server"
$ws.Range("C10").Value = "-"
$ws.Range("D10").Value = "Server"
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "-"
$ws.Range("G10").Value = "-"
$ws.Range("H10").Value = "Server"
$ws.Range("I10").Value = "-"
$ws.Range("J10").Value = "-"
$ws.Range("K10").Value = "-"
$ws.Range("L10").Value = "-"
$ws.Range("M10").Value = "-"
$ws.Range("N10").Value = "-"
$ws.Range("O10").Value = "-"
$ws.Range("B11").Value = "Server"
$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "-"
$ws.Range("G11").Value = "Server"
$ws.Range("H11").Value = "Server"
$ws.Range("I11").Value = "-"
$ws.Range("J11").Value = "-"
$ws.Range("K11").Value = "-"
$ws.Range("L11").Value = "-"
$ws.Range("M11").Value = "-"
$ws.Range("N11").Value = "-"
$ws.Range("O11").Value = "-"
$ws.Range("B12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "Server"
$ws.Range("F12").Value = "Server"
$ws.Range("H12").Value = "-"
$ws.Range("I12").Value = "Server"
$ws.Range("J12").Value = "-"
$ws.Range("L12").Value = "-"
$ws.Range("M12").Value = "-"
$ws.Range("N12").Value = "-"
$ws.Range("B13").Value = "-"
$ws.Range("C13").Value = "-"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = "-"
$ws.Range("F13").Value = "-"
$ws.Range("G13").Value = "-"
$ws.Range("H13").Value = "-"
$ws.Range("I13").Value = "-"
$ws.Range("J13").Value = "Server"
$ws.Range("K13").Value = "-"
$ws.Range("L13").Value = "-"
$ws.Range("M13").Value = "-"
$ws.Range("N13").Value = "-"
$ws.Range("O13").Value = "-"
$ws.Range("C14").Value = "Server"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"
$ws.Range("G14").Value = "-"
$ws.Range("H14").Value = "-"
$ws.Range("I14").Value = "Server"
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = "-"
$ws.Range("L14").Value = "-"
$ws.Range("M14").Value = "Server"
$ws.Range("N14").Value = "Server"
$ws.Range("O14").Value = "-"
$ws.Range("B15").Value = "Server"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"
$ws.Range("I15").Value = "-"
$ws.Range("J15").Value = "-"
$ws.Range("K15").Value = "Server"
$ws.Range("M15").Value = "-"
$ws.Range("N15").Value = "-"
$ws.Range("O15").Value = "Server"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"
$ws.Range("G16").Value = "-"
$ws.Range("H16").Value = "-"
$ws.Range("I16").Value = "-"
$ws.Range("J16").Value = "-"
$ws.Range("K16").Value = "Server"
$ws.Range("L16").Value = "Server"
$ws.Range("M16").Value = "-"
$ws.Range("N16").Value = "Server"
$ws.Range("O16").Value = "Server"
$ws.Range("B17").Value = "-"
$ws.Range("C17").Value = "-"
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = "-"
$ws.Range("F17").Value = "-"
$ws.Range("G17").Value = "-"
$ws.Range("H17").Value = "-"
$ws.Range("I17").Value = "-"
$ws.Range("J17").Value = "-"
$ws.Range("K17").Value = "-"
$ws.Range("L17").Value = "-"
$ws.Range("M17").Value = "Server"
$ws.Range("N17").Value = "-"
$ws.Range("O17").Value = "-"
$ws.Range("B18").Value = "Server"
$ws.Range("C18").Value = "Server"
$ws.Range("D18").Value = "Server"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "Server"
$ws.Range("H18").Value = "-"
$ws.Range("I18").Value = "Server"
$ws.Range("J18").Value = "-"
$ws.Range("K18").Value = "-"
$ws.Range("L18").Value = "-"
$ws.Range("M18").Value = "-"
$ws.Range("N18").Value = "-"
$ws.Range("O18").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("F19").Value = "-"
$ws.Range("G19").Value = "-"
$ws.Range("H19").Value = "-"
$ws.Range("I19").Value = "-"
$ws.Range("J19").Value = "Server"
$ws.Range("K19").Value = "-"
$ws.Range("L19").Value = "-"
$ws.Range("M19").Value = "-"
$ws.Range("N19").Value = "-"
$ws.Range("O19").Value = "-"
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"
$ws.Range("G20").Value = "-"
$ws.Range("H20").Value = "-"
$ws.Range("I20").Value = "-"
$ws.Range("J20").Value = "-"
$ws.Range("K20").Value = "-"
$ws.Range("L20").Value = "-"
$ws.Range("M20").Value = "-"
$ws.Range("N20").Value = "-"
$ws.Range("O20").Value = "-"
$ws.Range("B21").Value = "Server"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "Server"
$ws.Range("G21").Value = "Server"
$ws.Range("H21").Value = "-"
$ws.Range("I21").Value = "Server"
$ws.Range("J21").Value = "Server"
$ws.Range("L21").Value = "Expo"
$ws.Range("N21").Value = "Server"
$ws.Range("O21").Value = "Expo"
$ws.Range("B22").Value = "Server"
$ws.Range("C22").Value = "-"
$ws.Range("E22").Value = "Server"
$ws.Range("F22").Value = "-"
$ws.Range("G22").Value = "-"
$ws.Range("H22").Value = "-"
$ws.Range("I22").Value = "-"
$ws.Range("J22").Value = "-"
$ws.Range("K22").Value = "Server"
$ws.Range("L22").Value = "Server"
$ws.Range("M22").Value = "-"
$ws.Range("N22").Value = "-"
$ws.Range("O22").Value = "-"
$ws.Range("B23").Value = "Server"
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "Server"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "Server"
$ws.Range("I23").Value = "-"
$ws.Range("J23").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("L23").Value = "-"
$ws.Range("N23").Value = "-"
$ws.Range("O23").Value = "-"
$ws.Range("B24").Value = "-"
$ws.Range("C24").Value = "-"
$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = "-"
$ws.Range("I24").Value = "-"
$ws.Range("J24").Value = "-"
$ws.Range("K24").Value = "Server"
$ws.Range("L24").Value = "-"
$ws.Range("M24").Value = "-"
$ws.Range("N24").Value = "-"
$ws.Range("O24").Value = "Server"
$ws.Range("B25").Value = "-"
$ws.Range("C25").Value = "-"
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = "-"
$ws.Range("F25").Value = "-"
$ws.Range("G25").Value = "-"
$ws.Range("H25").Value = "Server"
$ws.Range("I25").Value = "-"
$ws.Range("J25").Value = "Server"
$ws.Range("K25").Value = "-"
$ws.Range("L25").Value = "Server"
$ws.Range("N25").Value = "-"
$ws.Range("O25").Value = "-"
$ws.Range("B26").Value = "H/G"
$ws.Range("C26").Value = "H/G"
$ws.Range("D26").Value = "-"
$ws.Range("F26").Value = "-"
$ws.Range("G26").Value = "-"
$ws.Range("H26").Value = "-"
$ws.Range("I26").Value = "-"
$ws.Range("K26").Value = "H/G"
$ws.Range("L26").Value = "-"
$ws.Range("M26").Value = "-"
$ws.Range("N26").Value = "-"
$ws.Range("O26").Value = "-"
$ws.Range("B27").Value = "-"
$ws.Range("C27").Value = "-"
$ws.Range("D27").Value = "-"
$ws.Range("E27").Value = "-"
$ws.Range("F27").Value = "-"
$ws.Range("G27").Value = "-"
$ws.Range("H27").Value = "-"
$ws.Range("I27").Value = "-"
$ws.Range("J27").Value = "-"
$ws.Range("K27").Value = "-"
$ws.Range("L27").Value = "-"
$ws.Range("M27").Value = "H/G"
$ws.Range("N27").Value = "-"
$ws.Range("O27").Value = "H/G"
$ws.Range("B28").Value = "-"
$ws.Range("C28").Value = "-"
$ws.Range("D28").Value = "-"
$ws.Range("E28").Value = "-"
$ws.Range("F28").Value = "-"
$ws.Range("G28").Value = "-"
$ws.Range("H28").Value = "H/G"
$ws.Range("I28").Value = "-"
$ws.Range("J28").Value = "-"
$ws.Range("K28").Value = "-"
$ws.Range("L28").Value = "H/G"
$ws.Range("M28").Value = "-"
$ws.Range("N28").Value = "-"
$ws.Range("O28").Value = "-"
$ws.Range("B29").Value = "-"
$ws.Range("C29").Value = "Runner"
$ws.Range("D29").Value = "-"
$ws.Range("E29").Value = "-"
$ws.Range("F29").Value = "-"
$ws.Range("G29").Value = "-"
$ws.Range("H29").Value = "-"
$ws.Range("I29").Value = "-"
$ws.Range("J29").Value = "-"
$ws.Range("K29").Value = "-"
$ws.Range("L29").Value = "-"
$ws.Range("M29").Value = "-"
$ws.Range("N29").Value = "-"
$ws.Range("O29").Value = "-"
$ws.Range("B30").Value = "-"
$ws.Range("C30").Value = "-"
$ws.Range("D30").Value = "H/G"
$ws.Range("E30").Value = "-"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "H/G"
$ws.Range("H30").Value = "-"
$ws.Range("I30").Value = "-"
$ws.Range("J30").Value = "-"
$ws.Range("K30").Value = "-"
$ws.Range("L30").Value = "-"
$ws.Range("M30").Value = "Runner"
$ws.Range("N30").Value = "-"
$ws.Range("B31").Value = "Runner"
$ws.Range("C31").Value = "-"
$ws.Range("D31").Value = "-"
$ws.Range("E31").Value = "Runner"
$ws.Range("F31").Value = "H/G"
$ws.Range("G31").Value = "-"
$ws.Range("H31").Value = "-"
$ws.Range("I31").Value = "H/G"
$ws.Range("J31").Value = "-"
$ws.Range("K31").Value = "-"
$ws.Range("L31").Value = "-"
$ws.Range("M31").Value = "-"
$ws.Range("N31").Value = "H/G"
$ws.Range("O31").Value = "-"
$ws.Range("B32").Value = "Expo"
$ws.Range("C32").Value = "Expo"
$ws.Range("D32").Value = "-"
$ws.Range("E32").Value = "-"
$ws.Range("F32").Value = "-"
$ws.Range("G32").Value = "-"
$ws.Range("H32").Value = "Expo"
$ws.Range("I32").Value = "-"
$ws.Range("K32").Value = "-"
$ws.Range("L32").Value = "-"
$ws.Range("M32").Value = "-"
$ws.Range("N32").Value = "-"
$ws.Range("O32").Value = "-"
